$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Measurement of Moment of Inertia" section ---
# Cells are written in the same order the original author entered them so
# the resulting shared-strings table lines up with the canonical file.

# Header row for the 40ms window table
$ws.Range("Q24").Value = "40ms window"

# t=0 row
$ws.Range("Q25").Value = "t=0"
$ws.Range("R25").Value = 9.4

# t=40ms row
$ws.Range("Q26").Value = "t=40ms"
$ws.Range("R26").Value = 4

# rms current
$ws.Range("Q28").Value = "rms current"
$ws.Range("R28").Value = 2.8

# LCR max
$ws.Range("Q31").Value = "LCR max"
$ws.Range("R31").Value = "2.7mH"

# min
$ws.Range("Q32").Value = "min"
$ws.Range("R32").Value = 2.3

# Section title
$ws.Range("Q21").Value = "Measurement of Moment of Inertia"

# Torque
$ws.Range("Q35").Value = "Torque "
$ws.Range("S35").Value = "N*m"

# Remaining header cells for the 40ms window table
$ws.Range("R24").Value = "ms"
$ws.Range("S24").Value = "2*ms"
$ws.Range("T24").Value = "rotation/s"
$ws.Range("U24").Value = "rad/s"

# Alpha
$ws.Range("Q36").Value = "Alpha"

# J
$ws.Range("Q38").Value = "J"

# Alpha units
$ws.Range("S36").Value = "rad/s^2"

# --- Formulas ---
$ws.Range("S25").Formula = "=R25*2"
$ws.Range("T25").Formula = "=1/S25*1000"
$ws.Range("U25").Formula = "=2*PI()*T25"

$ws.Range("S26").Formula = "=R26*2"
$ws.Range("T26").Formula = "=1/S26*1000"
$ws.Range("U26").Formula = "=2*PI()*T26"

$ws.Range("R35").Formula = "=R28/R3"
$ws.Range("R36").Formula = "=(U26-U25)/0.04"
$ws.Range("R38").Formula = "=R35/R36"

# --- View state changes ---
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 16520
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("S38").Select()
